$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 15
$ws.Range("B4").Value = 45
$ws.Range("B5").Value = 30

$ws.Range("B2:B5").HorizontalAlignment = -4108

$ws.Range("B5").Select()
